$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh snapshot. Values are written with a leading
# apostrophe to force text entry (so plain-looking decimals like '586.46'
# are not reinterpreted as numbers and lose their exact formatting), then
# the cell style is reset to Normal so no stray number-format style sticks.
$ws.Range("D2").Value = "'64.229.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.23%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.488.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.06%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'586.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'134.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.07%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.47%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.15%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.44%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.385"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.74%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.085.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.80%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000182"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.47%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.92%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.489.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.30%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'64.309.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.37%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'25.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -7.02%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'9.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.18%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'5.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.31%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.61%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'394.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.30%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -1.16%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.630.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.81%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'74.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.94%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.04%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.38%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.09%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.48%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -5.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.20%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'8.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.77%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.512.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.51%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E34").Value = "'  +0.10%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'23.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.75%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -4.28%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'6.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.19%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.15%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'166.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.80%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0780"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.27%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.805"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.35%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.12%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'25.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.39%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.83%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +1.90%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -3.73%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.457.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.46%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.07%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.55%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0260"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.32%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.73%  "
$ws.Range("E51").Style = "Normal"
